$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "卡鲁提拉号"
$ws.Range("B7").Font.Name = "宋体"
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = "2_42_370;1_225_340;3_97_266;3_177_266;3_257_266;4_337_266"
$ws.Range("F7").Value = "0;0;0;0;1;0;0"
